# Update Table 5 (Quantify variance explained by urbanization, Urbanization Score)
# values in the "flowertime_2022" ranova table, per re-run with 1000 iterations.

$d = $word.ActiveDocument

# 3.439 -> 3.443  (Block chi-sq)
$d.Content.Find.Execute("3.439", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3.443", 2)

# 0.329 -> 0.328  (Block p)
$d.Content.Find.Execute("0.329", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.328", 2)

# 0.034 -> 0.035  (Urbanization Score chi-sq)
$d.Content.Find.Execute("0.034", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.035", 2)

# 0.853 -> 0.851  (Urbanization Score p)
$d.Content.Find.Execute("0.853", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.851", 2)
